# Apply updated crypto market data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.614.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4741"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.92%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2904"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06483"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.01"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7386"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.01%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.33%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.872.85"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.182"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.600.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.96%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.36%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007487"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.119.43"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.234"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.203"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.32%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.90%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.899"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09922"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.346"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.510"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.243"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.089"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04772"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.44%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6928"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.38%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01855"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.758"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.226"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.26"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.971"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4159"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8343"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.61%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.385"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.65%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.971"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.66%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "918.42"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05669"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.71%  "
